$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $row, $date, $timestamp, $hour, $location, $value, $status)

    # Column A holds a date-like string ("2026-02-01") that must stay plain
    # text (matching the rest of the log), not get auto-converted to a real
    # Excel date serial. Force the cell to Text format before writing it.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $date

    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $status
}

# --- Proximity sheet: append rows 40-43 ---
$proximity = $wb.Worksheets.Item("Proximity")

Add-LogRow $proximity 40 "2026-02-01" "14:44:05" "14:00" "Living Room Main Door" "EXIT"  "User EXITED Living Room Main Door"
Add-LogRow $proximity 41 "2026-02-01" "14:44:10" "14:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $proximity 42 "2026-02-01" "14:44:20" "14:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $proximity 43 "2026-02-01" "14:44:23" "14:00" "Living Room Main Door" "EXIT"  "User EXITED Living Room Main Door"

# --- Camera sheet: append rows 26-29 ---
$camera = $wb.Worksheets.Item("Camera")

Add-LogRow $camera 26 "2026-02-01" "14:44:05" "14:00" "Living Room Main Door" "Image Captured" "Active"
Add-LogRow $camera 27 "2026-02-01" "14:44:07" "14:00" "Living Room Main Door" "Image Captured" "Active"
Add-LogRow $camera 28 "2026-02-01" "14:44:10" "14:00" "Living Room Main Door" "Image Received" "Active"
Add-LogRow $camera 29 "2026-02-01" "14:44:23" "14:00" "Living Room Main Door" "Image Captured" "Active"
